$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new keyword row (row 7), reusing the same keyword / appID pair
# already present in row 5 ("helix jump" / "com.singleton.helix").
$ws.Range("A7").Value = "helix jump"
$ws.Range("B7").Value = "com.singleton.helix"

# Move / extend the selection onto the newly added row, like Excel would
# leave it after typing the new values in.
$ws.Range("A7:B7").Select()
